$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-401) holds the "Förändrad" (changed) date as a date
# serial number. The commit updates every existing value of 45203
# (2023-10-04) to 45205 (2023-10-06).
for ($r = 2; $r -le 401; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45205
    }
}
